$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin price/volume figures are stored as plain text in this sheet (some
# prices even use two decimal points, e.g. "29.143.94"). Cells whose new
# text would otherwise be auto-parsed into a Double by Excel (losing the
# exact printed digits, e.g. "2.730" -> 2.73) are first marked as Text so
# the literal string is preserved.
$ws.Range("D2").Value = "29.143.94"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "1.831.95"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9991"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.37"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6640"
$ws.Range("E6").Value = "  -2.41%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9999"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07409"
$ws.Range("E8").Value = "  -0.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2935"
$ws.Range("E9").Value = "  -1.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.75"
$ws.Range("E10").Value = "  -1.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07741"
$ws.Range("E11").Value = "  +1.23%  "
$ws.Range("D12").Value = "1.834.71"
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.981"
$ws.Range("E13").Value = "  -0.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6676"
$ws.Range("E14").Value = "  -1.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.76"
$ws.Range("E15").Value = "  -4.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.091"
$ws.Range("E16").Value = "  -1.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008341"
$ws.Range("E17").Value = "  +1.60%  "
$ws.Range("D18").Value = "29.143.27"
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("D19").Value = "2.088.19"
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "227.36"
$ws.Range("E20").Value = "  -0.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.45"
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.153"
$ws.Range("E23").Value = "  -2.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9996"
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.49"
$ws.Range("E25").Value = "  -0.99%  "
$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1401"
$ws.Range("E26").Value = "  -2.06%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.594"
$ws.Range("E27").Value = "  -1.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.94"
$ws.Range("E28").Value = "  -0.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.506"
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.106"
$ws.Range("E30").Value = "  -3.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.031"
$ws.Range("E31").Value = "  -2.65%  "
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05310"
$ws.Range("E33").Value = "  -0.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.867"
$ws.Range("E34").Value = "  +1.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7486"
$ws.Range("E35").Value = "  -0.74%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.136"
$ws.Range("E36").Value = "  +0.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.644"
$ws.Range("E37").Value = "  -1.44%  "
$ws.Range("D38").Value = "1.272.49"
$ws.Range("E38").Value = "  -3.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01796"
$ws.Range("E39").Value = "  -1.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.730"
$ws.Range("E40").Value = "  +0.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9275"
$ws.Range("E41").Value = "  -0.85%  "
$ws.Range("B42").Value = "XinFinNetwork"
$ws.Range("C42").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.08557"
$ws.Range("E42").Value = "  +3.23%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.907"
$ws.Range("E43").Value = "  -2.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9996"
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.94"
$ws.Range("E45").Value = "  -2.86%  "
$ws.Range("D46").Value = "1.986.02"
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("E47").Value = "  -0.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.761"
$ws.Range("E48").Value = "  -0.38%  "
$ws.Range("E49").Value = "  -1.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "62.99"
$ws.Range("E50").Value = "  -1.76%  "
$ws.Range("E51").Value = "  -0.90%  "
